$d = $word.ActiveDocument

# The original doc has 5 paragraphs (title line, blank, two numbered
# bullet lines in Arial Narrow, trailing blank). The new doc collapses
# all of that down to 3 plain paragraphs: a bold title and two normal
# body lines. Rather than trying to strip numbering/fonts property by
# property, clear the whole body (keeping just the final paragraph
# mark) and retype fresh, unformatted paragraphs.

$d.Range(0, $d.Content.End - 1).Delete()

$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "Divided"
$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = "Stand together or don’t stand at all."
$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "(Can’t say I’m happy with it. Still working on it.)"

# Bold just the title paragraph, applied last so it doesn't leak into
# the paragraphs inserted afterwards.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Font.Bold = 1
